# Commit: "changing small b to B"
# The slide's "observed effect (b*)" label is re-labelled to use an
# upper-case B, matching the adjoining between-study-variability formula's
# notation (B* for the observed effect).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -like "observed effect*") {
        $target = $shp
        break
    }
}

if ($target -eq $null) {
    $target = $s.Shapes.Item("Tekstvak 7")
}

$target.TextFrame.TextRange.Text = "observed effect (B*)"
